# 1. Refresh the "generated at" timestamp in the footer.
$d = $word.ActiveDocument
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute("2025-06-30 12:12Z / ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "2025-07-02 02:48Z / ", 2)

# 2. Add the standard PubMed "b / i / sub / sup / u" inline character styles
#    (regression test for round-tripping these helper styles).
$defs = @(
    @{ Id = "b";   Bold = $true },
    @{ Id = "i";   Italic = $true },
    @{ Id = "sub"; Subscript = $true },
    @{ Id = "sup"; Superscript = $true },
    @{ Id = "u";   Underline = $true }
)

foreach ($def in $defs) {
    $style = $d.Styles.Add($def.Id, 2)
    $style.BaseStyle = $d.Styles("DefaultParagraphFont")
    $style.Priority = 1
    $style.QuickStyle = $true

    if ($def.Bold) { $style.Font.Bold = $true }
    if ($def.Italic) { $style.Font.Italic = $true }
    if ($def.Subscript) { $style.Font.Subscript = $true }
    if ($def.Superscript) { $style.Font.Superscript = $true }
    if ($def.Underline) { $style.Font.Underline = 1 }
}
